$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("T2").Value = "maa://22742 (90.75), *maa://20791 (62.82)"
$ws.Range("D3").Value = "maa://40192 (96.83), maa://36987 (96.15), maa://39849 (88.89)"
$ws.Range("T3").Value = "maa://24617 (90.0), **maa://20790 (43.48), ***maa://37170 (16.18), maa://45854 (85.71)"
$ws.Range("T4").Value = "maa://32509 (95.8), maa://27295 (86.67), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)"
$ws.Range("X4").Value = "**maa://32495 (48.71), ***maa://31785 (22.22), maa://43217 (91.14), ***maa://36683 (28.26)"
$ws.Range("D5").Value = "maa://21245 (84.55), maa://22744 (84.62)"
$ws.Range("A8").Value = "更新日期：2025.04.04 13:19:56"
$ws.Range("P8").Value = "maa://32931 (84.17), *maa://21916 (62.12), maa://23252 (91.18), maa://37496 (96.97), **maa://22759 (45.45)"
$ws.Range("X8").Value = "maa://21411 (96.01)"
$ws.Range("T10").Value = "maa://27395 (96.6), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("T11").Value = "maa://22747 (92.64), maa://22501 (97.85), maa://45521 (86.36)"
$ws.Range("X12").Value = "maa://22753 (91.16), *maa://21485 (75.35), maa://37962 (90.7)"
$ws.Range("D13").Value = "maa://24999 (92.19), maa://36673 (92.41), maa://25001 (85.92)"
$ws.Range("H13").Value = "*maa://21248 (73.6), **maa://22728 (46.67)"
$ws.Range("D15").Value = "*maa://22743 (78.24), maa://22734 (84.17), *maa://30808 (64.18), *maa://36048 (51.47), maa://45058 (93.33)"
$ws.Range("AF15").Value = "maa://21364 (80.99), *maa://36666 (77.12), *maa://22766 (68.33)"
$ws.Range("P16").Value = "maa://28504 (91.67)"
$ws.Range("AB16").Value = "maa://26228 (95.19)"
$ws.Range("L18").Value = "maa://22466 (90.96), *maa://22732 (51.55)"
$ws.Range("AB19").Value = "*maa://30709 (66.22), *maa://36668 (58.02)"
$ws.Range("D20").Value = "maa://21432 (90.66), maa://25198 (93.75), *maa://20795 (50.77), maa://36680 (91.18)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (91.3)"
$ws.Range("AB21").Value = "maa://21443 (81.56), ***maa://23820 (30.0)"
$ws.Range("L23").Value = "maa://39756 (95.68), maa://39875 (94.59)"
$ws.Range("X24").Value = "maa://29988 (83.71), maa://23504 (93.39), **maa://22892 (40.54), *maa://25141 (77.1), *maa://36663 (77.5), ***maa://22815 (23.08)"
$ws.Range("H25").Value = "*maa://29063 (72.94), *maa://25311 (74.77), ***maa://22725 (4.76), *maa://45047 (66.67)"
$ws.Range("D26").Value = "maa://41802 (91.3)"
$ws.Range("AF29").Value = "*maa://24080 (68.93), maa://42865 (81.94), ***maa://34960 (8.33)"
$ws.Range("H32").Value = "maa://21895 (97.12), maa://36667 (97.73), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("P34").Value = "maa://48817 (88.24)"
$ws.Range("L37").Value = "maa://45718 (97.67), *maa://47069 (75.0), maa://45789 (100.0)"
$ws.Range("H53").Value = "maa://32534 (94.27), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.16)"
$ws.Range("H59").Value = "maa://31270 (94.81), maa://27746 (82.46)"
$ws.Range("H60").Value = "*maa://40438 (70.77)"
